# -----------------------------------------------------------------------
# screen_shot.docx edit:
#   * Insert a blank "noProof" paragraph above each picture.
#   * Crop / resize the two pictures (srcRect + new extents) and give the
#     picture shapes the "no line / shadow-obscured" formatting that Word
#     stamps on when a picture style is touched in the UI.
#   * Move the "_GoBack" bookmark from the end of Picture 2's paragraph to
#     the start of that paragraph.
# -----------------------------------------------------------------------

$d = $word.ActiveDocument

# --- Discover the identifiers already present in the document instead of
#     hard-coding them, so the script keeps working if the relationship
#     ids / shape ids / anchor ids ever differ from this particular file.
$srcXml = $d.Content.WordOpenXML

$ro = [System.Text.RegularExpressions.RegexOptions]::Singleline

$drawingPattern = '<wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="([0-9A-F]+)" wp14:editId="([0-9A-F]+)">.*?<wp:docPr id="(\d+)" name="([^"]*)"/>.*?<a:blip r:embed="(rId\d+)"/>'
$drawingMatches = [regex]::Matches($srcXml, $drawingPattern, $ro)

$pic1AnchorId = $drawingMatches[0].Groups[1].Value
$pic1EditId   = $drawingMatches[0].Groups[2].Value
$pic1DocPrId  = $drawingMatches[0].Groups[3].Value
$pic1DocPrNm  = $drawingMatches[0].Groups[4].Value
$pic1RId      = $drawingMatches[0].Groups[5].Value

$pic2AnchorId = $drawingMatches[1].Groups[1].Value
$pic2EditId   = $drawingMatches[1].Groups[2].Value
$pic2DocPrId  = $drawingMatches[1].Groups[3].Value
$pic2DocPrNm  = $drawingMatches[1].Groups[4].Value
$pic2RId      = $drawingMatches[1].Groups[5].Value

if ($srcXml -match '<w:document[^>]*>') {
    $docElementOpenTag = $matches[0]
} else {
    $docElementOpenTag = '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
}

# --- Build the replacement body (everything up to, but not including,
#     the trailing sectPr - $d.Content never includes the final section
#     mark, so the existing sectPr is left completely untouched).
$newBodyXml = @"
$docElementOpenTag
  <w:body>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <w:drawing>
          <wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="$pic1AnchorId" wp14:editId="$pic1EditId">
            <wp:extent cx="5669280" cy="2537460"/>
            <wp:effectExtent l="0" t="0" r="7620" b="0"/>
            <wp:docPr id="$pic1DocPrId" name="$pic1DocPrNm"/>
            <wp:cNvGraphicFramePr>
              <a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/>
            </wp:cNvGraphicFramePr>
            <a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main">
              <a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture">
                <pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture">
                  <pic:nvPicPr>
                    <pic:cNvPr id="1" name=""/>
                    <pic:cNvPicPr/>
                  </pic:nvPicPr>
                  <pic:blipFill rotWithShape="1">
                    <a:blip r:embed="$pic1RId"/>
                    <a:srcRect t="5015" r="4616" b="19089"/>
                    <a:stretch/>
                  </pic:blipFill>
                  <pic:spPr bwMode="auto">
                    <a:xfrm>
                      <a:off x="0" y="0"/>
                      <a:ext cx="5669280" cy="2537460"/>
                    </a:xfrm>
                    <a:prstGeom prst="rect">
                      <a:avLst/>
                    </a:prstGeom>
                    <a:ln>
                      <a:noFill/>
                    </a:ln>
                    <a:extLst>
                      <a:ext uri="{53640926-AAD7-44D8-BBD7-CCE9431645EC}">
                        <a14:shadowObscured xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main"/>
                      </a:ext>
                    </a:extLst>
                  </pic:spPr>
                </pic:pic>
              </a:graphicData>
            </a:graphic>
          </wp:inline>
        </w:drawing>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <w:drawing>
          <wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="$pic2AnchorId" wp14:editId="$pic2EditId">
            <wp:extent cx="5814060" cy="2796540"/>
            <wp:effectExtent l="0" t="0" r="0" b="3810"/>
            <wp:docPr id="$pic2DocPrId" name="$pic2DocPrNm"/>
            <wp:cNvGraphicFramePr>
              <a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/>
            </wp:cNvGraphicFramePr>
            <a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main">
              <a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture">
                <pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture">
                  <pic:nvPicPr>
                    <pic:cNvPr id="1" name=""/>
                    <pic:cNvPicPr/>
                  </pic:nvPicPr>
                  <pic:blipFill rotWithShape="1">
                    <a:blip r:embed="$pic2RId"/>
                    <a:srcRect l="-128" t="7066" r="2308" b="9288"/>
                    <a:stretch/>
                  </pic:blipFill>
                  <pic:spPr bwMode="auto">
                    <a:xfrm>
                      <a:off x="0" y="0"/>
                      <a:ext cx="5814060" cy="2796540"/>
                    </a:xfrm>
                    <a:prstGeom prst="rect">
                      <a:avLst/>
                    </a:prstGeom>
                    <a:ln>
                      <a:noFill/>
                    </a:ln>
                    <a:extLst>
                      <a:ext uri="{53640926-AAD7-44D8-BBD7-CCE9431645EC}">
                        <a14:shadowObscured xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main"/>
                      </a:ext>
                    </a:extLst>
                  </pic:spPr>
                </pic:pic>
              </a:graphicData>
            </a:graphic>
          </wp:inline>
        </w:drawing>
      </w:r>
    </w:p>
  </w:body>
</w:document>
"@

$pkgXml = @"
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>$newBodyXml</pkg:xmlData></pkg:part></pkg:package>
"@

$null = $d.Content.InsertXML($pkgXml)
